$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty Epochs value for the row 8 "Adam / 12288" run
$ws.Range("C8").Value = 10

# Insert a new results row (the "Adam / 12288 / Epochs 20" run) right after
# row 8, pushing the SGD row and the explanatory note block down by one row.
$ws.Rows("9:9").Insert()

# Copy row 8's formatting onto the freshly-inserted row 9 so the new row
# keeps the table's border/number formatting instead of Excel's blank
# default for inserted rows.
$ws.Range("A8:H8").Copy()
$ws.Range("A9:H9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row's data.
$ws.Range("A9").Value = "CodeGPTPy"
$ws.Range("B9").Value = 12288
$ws.Range("C9").Value = 20
$ws.Range("D9").Value = "Adam"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0.667
$ws.Range("G9").Value = 0.653
$ws.Range("H9").Formula = "=F9-G9"

# Match the saved selection from the source workbook.
[void]$ws.Range("H9").Select()
